# "organization and predict-observe updates"
# Remove the ten "Straight Connector 16".."Straight Connector 25" straight-line
# connector shapes from the second slide (they formed a row guide/grid that is
# no longer needed). All other shapes on the slide are left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$connectorNames = @(
    "Straight Connector 16",
    "Straight Connector 17",
    "Straight Connector 18",
    "Straight Connector 19",
    "Straight Connector 20",
    "Straight Connector 21",
    "Straight Connector 22",
    "Straight Connector 23",
    "Straight Connector 24",
    "Straight Connector 25"
)

foreach ($name in $connectorNames) {
    $s.Shapes.Item($name).Delete()
}
